$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C12").Value = "bicycle"
$ws.Range("B12").Select() | Out-Null
